$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Amount values for the first two data rows were zeroed out (BDD test fixture data).
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 0

# Re-apply the "Normal" style to the SL No column (A1:A22) -- this is what
# produced the extra cellXfs entry in the original commit (same visual
# appearance, just written out as an explicit style again).
$ws.Range("A1:A22").Style = "Normal"

# Cursor ended up on E3 when the fixture was saved.
$ws.Range("E3").Select()
